# Auto-generated edit script: updates Leve profit-calc values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC!row 64
$ws_ALC.Range("H64").Value = 4664.8887
$ws_ALC.Range("I64").Value = 4748
$ws_ALC.Range("J64").Value = 4000
$ws_ALC.Range("K64").Value = 4748
$ws_ALC.Range("L64").Value = 4000
$ws_ALC.Range("M64").Value = -4500
$ws_ALC.Range("N64").Value = -4496

# ALC!row 67
$ws_ALC.Range("H67").Value = 4664.8887
$ws_ALC.Range("I67").Value = 4748
$ws_ALC.Range("J67").Value = 4000
$ws_ALC.Range("K67").Value = 4748
$ws_ALC.Range("L67").Value = 4000
$ws_ALC.Range("M67").Value = -3890
$ws_ALC.Range("N67").Value = -5716

# ALC!row 74
$ws_ALC.Range("H74").Value = 3962.5
$ws_ALC.Range("I74").Value = 3957.1428
$ws_ALC.Range("J74").Value = 4000
$ws_ALC.Range("K74").Value = 3957.1428
$ws_ALC.Range("L74").Value = 4000
$ws_ALC.Range("M74").Value = -3021.1428
$ws_ALC.Range("N74").Value = -5872

# ALC!row 76
$ws_ALC.Range("H76").Value = 3150.6191
$ws_ALC.Range("I76").Value = 3101.6667
$ws_ALC.Range("K76").Value = 3101.6667
$ws_ALC.Range("M76").Value = -2786.6667

# ALC!row 77
$ws_ALC.Range("H77").Value = 3962.5
$ws_ALC.Range("I77").Value = 3957.1428
$ws_ALC.Range("J77").Value = 4000
$ws_ALC.Range("K77").Value = 19785.714
$ws_ALC.Range("L77").Value = 20000
$ws_ALC.Range("M77").Value = -15105.714
$ws_ALC.Range("N77").Value = -29360

# ALC!row 79
$ws_ALC.Range("H79").Value = 3150.6191
$ws_ALC.Range("I79").Value = 3101.6667
$ws_ALC.Range("K79").Value = 3101.6667
$ws_ALC.Range("M79").Value = -2009.6667

# ALC!row 112
$ws_ALC.Range("H112").Value = 43479612
$ws_ALC.Range("J112").Value = 50001468
$ws_ALC.Range("L112").Value = 150004404
$ws_ALC.Range("N112").Value = -150006620

# ALC!row 113
$ws_ALC.Range("H113").Value = 2220
$ws_ALC.Range("I113").Value = 3000
$ws_ALC.Range("J113").Value = 2025
$ws_ALC.Range("K113").Value = 3000
$ws_ALC.Range("L113").Value = 2025
$ws_ALC.Range("M113").Value = 254
$ws_ALC.Range("N113").Value = -8533

# ALC!row 138
$ws_ALC.Range("H138").Value = 2750.7793
$ws_ALC.Range("J138").Value = 2820.3428
$ws_ALC.Range("L138").Value = 8461.028399999999
$ws_ALC.Range("N138").Value = -18741.0284

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM!row 61
$ws_ARM.Range("H61").Value = 3282.4634
$ws_ARM.Range("I61").Value = 2881.238
$ws_ARM.Range("K61").Value = 2881.238
$ws_ARM.Range("M61").Value = -2669.238

# ARM!row 117
$ws_ARM.Range("H117").Value = 20045
$ws_ARM.Range("J117").Value = 20045
$ws_ARM.Range("L117").Value = 20045
$ws_ARM.Range("N117").Value = -29223

# ARM!row 136
$ws_ARM.Range("H136").Value = 3282.4634
$ws_ARM.Range("I136").Value = 2881.238
$ws_ARM.Range("K136").Value = 8643.714
$ws_ARM.Range("M136").Value = -6093.714

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM!row 86
$ws_BSM.Range("H86").Value = 1692.8572
$ws_BSM.Range("I86").Value = 1600
$ws_BSM.Range("J86").Value = 2250
$ws_BSM.Range("K86").Value = 1600
$ws_BSM.Range("L86").Value = 2250
$ws_BSM.Range("M86").Value = -477
$ws_BSM.Range("N86").Value = -4496

# BSM!row 89
$ws_BSM.Range("H89").Value = 1692.8572
$ws_BSM.Range("I89").Value = 1600
$ws_BSM.Range("J89").Value = 2250
$ws_BSM.Range("K89").Value = 8000
$ws_BSM.Range("L89").Value = 11250
$ws_BSM.Range("M89").Value = -2384
$ws_BSM.Range("N89").Value = -22482

# BSM!row 105
$ws_BSM.Range("H105").Value = 1708.9
$ws_BSM.Range("I105").Value = 1598.1666
$ws_BSM.Range("J105").Value = 1875
$ws_BSM.Range("K105").Value = 1598.1666
$ws_BSM.Range("L105").Value = 1875
$ws_BSM.Range("M105").Value = 148.8334
$ws_BSM.Range("N105").Value = -5369

# BSM!row 134
$ws_BSM.Range("H134").Value = 69275.7
$ws_BSM.Range("I134").Value = 83538.07000000001
$ws_BSM.Range("J134").Value = 2718
$ws_BSM.Range("K134").Value = 250614.21
$ws_BSM.Range("L134").Value = 8154
$ws_BSM.Range("M134").Value = -248079.21
$ws_BSM.Range("N134").Value = -13224

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP!row 62
$ws_CRP.Range("H62").Value = 3246.8667
$ws_CRP.Range("I62").Value = 2970.3
$ws_CRP.Range("J62").Value = 3800
$ws_CRP.Range("K62").Value = 2970.3
$ws_CRP.Range("L62").Value = 3800
$ws_CRP.Range("M62").Value = -2346.3
$ws_CRP.Range("N62").Value = -5048

# CRP!row 65
$ws_CRP.Range("H65").Value = 3246.8667
$ws_CRP.Range("I65").Value = 2970.3
$ws_CRP.Range("J65").Value = 3800
$ws_CRP.Range("K65").Value = 14851.5
$ws_CRP.Range("L65").Value = 19000
$ws_CRP.Range("M65").Value = -11731.5
$ws_CRP.Range("N65").Value = -25240

# CRP!row 132
$ws_CRP.Range("H132").Value = 2007.62
$ws_CRP.Range("I132").Value = 1842.0476
$ws_CRP.Range("J132").Value = 2876.875
$ws_CRP.Range("K132").Value = 5526.142800000001
$ws_CRP.Range("L132").Value = 8630.625
$ws_CRP.Range("M132").Value = -2996.142800000001
$ws_CRP.Range("N132").Value = -13690.625

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL!row 4
$ws_CUL.Range("H4").Value = 893.4211
$ws_CUL.Range("I4").Value = 465
$ws_CUL.Range("J4").Value = 1007.6667
$ws_CUL.Range("K4").Value = 1395
$ws_CUL.Range("L4").Value = 3023.0001
$ws_CUL.Range("M4").Value = -1283
$ws_CUL.Range("N4").Value = -3247.0001

# CUL!row 75
$ws_CUL.Range("H75").Value = 1001.5
$ws_CUL.Range("I75").Value = 503
$ws_CUL.Range("J75").Value = 1500
$ws_CUL.Range("K75").Value = 1509
$ws_CUL.Range("L75").Value = 4500
$ws_CUL.Range("M75").Value = -511
$ws_CUL.Range("N75").Value = -6496

# CUL!row 78
$ws_CUL.Range("H78").Value = 1001.5
$ws_CUL.Range("I78").Value = 503
$ws_CUL.Range("J78").Value = 1500
$ws_CUL.Range("K78").Value = 4527
$ws_CUL.Range("L78").Value = 13500
$ws_CUL.Range("M78").Value = 465
$ws_CUL.Range("N78").Value = -23484

# CUL!row 131
$ws_CUL.Range("H131").Value = 2485.8594
$ws_CUL.Range("I131").Value = 15157.5
$ws_CUL.Range("J131").Value = 1641.0834
$ws_CUL.Range("K131").Value = 45472.5
$ws_CUL.Range("L131").Value = 4923.2502
$ws_CUL.Range("M131").Value = -40432.5
$ws_CUL.Range("N131").Value = -15003.2502

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM!row 70
$ws_GSM.Range("H70").Value = 4708.2812
$ws_GSM.Range("I70").Value = 4625.4287
$ws_GSM.Range("J70").Value = 4772.722
$ws_GSM.Range("K70").Value = 4625.4287
$ws_GSM.Range("L70").Value = 4772.722
$ws_GSM.Range("M70").Value = -4355.4287
$ws_GSM.Range("N70").Value = -5312.722

# GSM!row 73
$ws_GSM.Range("H73").Value = 4708.2812
$ws_GSM.Range("I73").Value = 4625.4287
$ws_GSM.Range("J73").Value = 4772.722
$ws_GSM.Range("K73").Value = 4625.4287
$ws_GSM.Range("L73").Value = 4772.722
$ws_GSM.Range("M73").Value = -3689.4287
$ws_GSM.Range("N73").Value = -6644.722

# GSM!row 80
$ws_GSM.Range("H80").Value = 86070.5
$ws_GSM.Range("I80").Value = 2530
$ws_GSM.Range("J80").Value = 253151.5
$ws_GSM.Range("K80").Value = 2530
$ws_GSM.Range("L80").Value = 253151.5
$ws_GSM.Range("M80").Value = -1532
$ws_GSM.Range("N80").Value = -255147.5

# GSM!row 83
$ws_GSM.Range("H83").Value = 86070.5
$ws_GSM.Range("I83").Value = 2530
$ws_GSM.Range("J83").Value = 253151.5
$ws_GSM.Range("K83").Value = 12650
$ws_GSM.Range("L83").Value = 1265757.5
$ws_GSM.Range("M83").Value = -7658
$ws_GSM.Range("N83").Value = -1275741.5
